# Fruta / hortaliza, semanal
# Insert a new data row above the current row 8 (pushing existing rows 8-17
# down to 9-18) and populate it with this week's observation for
# "Comercializadora del Agro de Limarí - Arándano (blue)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8:17 down one row, inserting a fresh (blank) row 8.
$ws.Rows(8).Insert()

# Populate the newly-inserted row 8 with the new weekly record.
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44895
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = "Berries"
$ws.Range("I8").Value = 100101001
$ws.Range("J8").Value = "Arándano (blue)"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 240
$ws.Range("N8").Value = 3000
$ws.Range("O8").Value = 3500
$ws.Range("P8").Value = 3250
$ws.Range("Q8").Value = "$/bandeja 2 kilos"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 1625
$ws.Range("T8").Value = 2
